$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at 311, shifting existing rows 311:325 down to 312:326
$ws.Rows.Item(311).Insert()

# Populate the newly inserted row 311 with its data
$ws.Cells.Item(311, 1).Value = 5
$ws.Cells.Item(311, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(311, 3).Value = 'Maule'
$ws.Cells.Item(311, 4).Value = 44585
$ws.Cells.Item(311, 5).Value = 7
$ws.Cells.Item(311, 6).Value = 'Fruta'
$ws.Cells.Item(311, 7).Value = 100109
$ws.Cells.Item(311, 8).Value = 'Uva'
$ws.Cells.Item(311, 9).Value = 100109001
$ws.Cells.Item(311, 10).Value = 'Uva'
$ws.Cells.Item(311, 11).Value = 'Superior Seedless'
$ws.Cells.Item(311, 12).Value = 'Primera'
$ws.Cells.Item(311, 13).Value = 100
$ws.Cells.Item(311, 14).Value = 12000
$ws.Cells.Item(311, 15).Value = 12000
$ws.Cells.Item(311, 16).Value = 12000
$ws.Cells.Item(311, 17).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(311, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(311, 19).Value = 667
$ws.Cells.Item(311, 20).Value = 18
